$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 507
$ws.Range("I2").Value = 1328
$ws.Range("J2").Value = 5601
$ws.Range("L2").Value = 1518
$ws.Range("M2").Value = 88
$ws.Range("N2").Value = 989
$ws.Range("O2").Value = 3
$ws.Range("P2").Value = 24
$ws.Range("Q2").Value = 9
$ws.Range("R2").Value = 64
$ws.Range("S2").Value = 624
$ws.Range("T2").Value = 979
$ws.Range("U2").Value = 75
$ws.Range("V2").Value = 8785
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 8719
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 138
$ws.Range("AA2").Value = 56
